# Update report header labels and remove the now-unused "N° Equipos:" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Fecha:"
$ws.Range("A5").Value = "Laboratorio:"
$ws.Range("A6").Value = "N° de Registros:"

# Remove the old row 7 ("N° Equipos:") entirely, shrinking the sheet to A1:A6.
$ws.Range("A7").EntireRow.Delete()

# Move the active selection to A4, matching the final workbook state.
$ws.Range("A4").Select()
